$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (row 2..115): replace sequential month index (1..114) with
# Excel date serials for the 1st of each month, Jan 2014 .. Jun 2023 ---
$dates = @(41640,41671,41699,41730,41760,41791,41821,41852,41883,41913,41944,41974,42005,42036,42064,42095,42125,42156,42186,42217,42248,42278,42309,42339,42370,42401,42430,42461,42491,42522,42552,42583,42614,42644,42675,42705,42736,42767,42795,42826,42856,42887,42917,42948,42979,43009,43040,43070,43101,43132,43160,43191,43221,43252,43282,43313,43344,43374,43405,43435,43466,43497,43525,43556,43586,43617,43647,43678,43709,43739,43770,43800,43831,43862,43891,43922,43952,43983,44013,44044,44075,44105,44136,44166,44197,44228,44256,44287,44317,44348,44378,44409,44440,44470,44501,44531,44562,44593,44621,44652,44682,44713,44743,44774,44805,44835,44866,44896,44927,44958,44986,45017,45047,45078)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $dates[$i]
}

# Apply a date display format ("m/d/yyyy"-style => built-in numFmtId 14)
# to the whole column (header included), matching the style used
# throughout column A.
$ws.Range("A1").NumberFormat = "mm-dd-yy"
$ws.Range("A1").Copy()
$ws.Range("A1:A115").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column width: column A grows to fit the new date strings ---
$ws.Columns("A").ColumnWidth = 25.43

# --- Sheet view: scroll position + active selection ---
$ws.Range("C22").Select()
